$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list (price + 1h volume-change columns) refreshed by the
# scheduled GitHub Actions scrape. Values are stored as text (matching the
# existing inline-string cells), so we force text entry with a leading
# apostrophe to stop Excel from auto-converting numeric-looking strings.

$ws.Range("D2").Value = "'330.48"
$ws.Range("E2").Value = "'0.41%"

$ws.Range("D3").Value = "'41.23"
$ws.Range("E3").Value = "'1.20%"

$ws.Range("D4").Value = "'5.679"
$ws.Range("E4").Value = "'-0.94%"

$ws.Range("D5").Value = "'0.08064"
$ws.Range("E5").Value = "'-0.68%"

$ws.Range("D6").Value = "'2.023"
$ws.Range("E6").Value = "'2.89%"

$ws.Range("D7").Value = "'8.741"
$ws.Range("E7").Value = "'-0.19%"

$ws.Range("D8").Value = "'4.526"
$ws.Range("E8").Value = "'-1.62%"

$ws.Range("D9").Value = "'2.932"
$ws.Range("E9").Value = "'-0.32%"

$ws.Range("D10").Value = "'0.9249"
$ws.Range("E10").Value = "'-2.12%"

$ws.Range("D11").Value = "'0.1264"
$ws.Range("E11").Value = "'-2.85%"

$ws.Range("D12").Value = "'0.1938"
$ws.Range("E12").Value = "'-2.71%"

$ws.Range("D13").Value = "'8.267"
$ws.Range("E13").Value = "'-7.32%"

$ws.Range("D14").Value = "'0.09295"
$ws.Range("E14").Value = "'-1.06%"

$ws.Range("D15").Value = "'0.03707"
$ws.Range("E15").Value = "'5.52%"

$ws.Range("E16").Value = "'9.36%"

$ws.Range("D17").Value = "'0.001311"
$ws.Range("E17").Value = "'-0.70%"

$ws.Range("D18").Value = "'0.006240"
$ws.Range("E18").Value = "'-0.24%"

$ws.Range("D19").Value = "'3.364"
$ws.Range("E19").Value = "'0.13%"

$ws.Range("D20").Value = "'0.3473"
$ws.Range("E20").Value = "'-2.56%"

$ws.Range("D21").Value = "'0.1419"
$ws.Range("E21").Value = "'-0.16%"

$ws.Range("D22").Value = "'0.2656"
$ws.Range("E22").Value = "'10.07%"

$ws.Range("D23").Value = "'0.04409"
$ws.Range("E23").Value = "'-0.60%"

$ws.Range("D24").Value = "'0.001262"
$ws.Range("E24").Value = "'0.04%"

$ws.Range("D25").Value = "'0.004339"
$ws.Range("E25").Value = "'-1.03%"

$ws.Range("D26").Value = "'0.0001241"
$ws.Range("E26").Value = "'13.75%"

$ws.Range("D39").Value = "'0.02844"
$ws.Range("E39").Value = "'15.43%"

$ws.Range("D40").Value = "'0.05459"
$ws.Range("E40").Value = "'2.93%"

$ws.Range("D41").Value = "'0.007761"

$ws.Range("D42").Value = "'0.009964"
$ws.Range("E42").Value = "'15.01%"

$ws.Range("D43").Value = "'0.1419"
$ws.Range("E43").Value = "'-1.10%"

$ws.Range("D44").Value = "'0.002236"
$ws.Range("E44").Value = "'9.05%"

$ws.Range("D45").Value = "'0.01189"
$ws.Range("E45").Value = "'13.06%"

$ws.Range("D46").Value = "'0.00006779"
$ws.Range("E46").Value = "'-1.68%"

$ws.Range("E47").Value = "'-0.01%"

$ws.Range("D48").Value = "'0.003005"
$ws.Range("E48").Value = "'-14.27%"

$ws.Range("D49").Value = "'0.002284"
$ws.Range("E49").Value = "'34.20%"

$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.01%"

$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.01%"
